$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A53").Value = " 03-11-20"
$ws.Range("A54").Value = " 06-11-20"
